$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for the "Toda a série histórica" block (rows 2-7), now reordered by value
# and re-based on 2012/2022, with the Categoria column (D) removed.
$ws.Range("A2").Value = "Atividades financeiras, de seguros e serviços relacionados"
$ws.Range("B2").Value = 35.19143837025082
$ws.Range("C2").Value = "2012 / 2022"

$ws.Range("A3").Value = "Atividades imobiliárias"
$ws.Range("B3").Value = 31.76377306549257
$ws.Range("C3").Value = "2012 / 2022"

$ws.Range("A4").Value = "Informação e comunicação"
$ws.Range("B4").Value = 21.77658698762782
$ws.Range("C4").Value = "2012 / 2022"

$ws.Range("A5").Value = "Agropecuária"
$ws.Range("B5").Value = 9.730306427073359
$ws.Range("C5").Value = "2012 / 2022"

$ws.Range("A6").Value = "Administração, defesa, educação e saúde públicas e seguridade social"
$ws.Range("B6").Value = 4.760483082368495
$ws.Range("C6").Value = "2012 / 2022"

$ws.Range("A7").Value = "Eletricidade e gás, água, esgoto, atividades de gestão de resíduos e descontaminação"
$ws.Range("B7").Value = 2.412539862254022
$ws.Range("C7").Value = "2012 / 2022"

# Remove column D entirely (Categoria header + all values)
$ws.Range("D1:D13").Delete()
